$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old A4 cell (row 4 disappears entirely once it's empty).
$ws.Range("A4").ClearContents()

# New column D: vertical-alignment-only styles, empty cells (value-less).
$ws.Range("D1").VerticalAlignment = -4160   # xlVAlignTop
$ws.Range("D2").VerticalAlignment = -4108   # xlVAlignCenter
$ws.Range("D3").WrapText = $False           # produces an "empty" alignment xf

# New column E: same text as the rest of the sheet ("abc"), reusing the
# D1/D2 styles for E1/E2, default style for E3.
$ws.Range("E1").Value = "abc"
$ws.Range("E1").VerticalAlignment = -4160   # xlVAlignTop (same as D1 -> style 10)
$ws.Range("E2").Value = "abc"
$ws.Range("E2").VerticalAlignment = -4108   # xlVAlignCenter (same as D2 -> style 11)
$ws.Range("E3").Value = "abc"

# Update the active selection to D1, matching the edited workbook.
$ws.Range("D1").Select() | Out-Null
